$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L")
for ($r=1; $r -le 16; $r++) {
  $line = "Row " + $r + ": "
  foreach ($col in $cols) {
    $c = $ws.Range($col + $r)
    $bl = $c.Borders.Item(9).LineStyle
    $bw = $c.Borders.Item(9).Weight
    $line += $col + "(" + $bl + "," + $bw + ") "
  }
  Write-Host $line
}
